$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number would be auto-converted
# to a numeric type by Excel on assignment, losing the original text formatting
# (e.g. "4.500" -> 4.5). These addresses get a temporary Text number format so
# the literal text is preserved, then the format is restored to the original style.
$forceTextAddrs = @("D5", "D6", "D8", "D9", "D11", "D12", "D15", "D20", "D21", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D43", "D46", "D48", "D49", "D50", "D51")

function Set-CellText($addr, $value) {
    $cell = $ws.Range($addr)
    if ($forceTextAddrs -contains $addr) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $value
    }
}

Set-CellText "D2" "26.059.24"
Set-CellText "E2" "  -0.53%  "
Set-CellText "D3" "1.651.53"
Set-CellText "E3" "  -0.44%  "
Set-CellText "E4" "  -0.23%  "
Set-CellText "D5" "217.24"
Set-CellText "E5" "  +0.05%  "
Set-CellText "D6" "0.5259"
Set-CellText "E6" "  +1.90%  "
Set-CellText "E7" "  -0.15%  "
Set-CellText "D8" "0.2597"
Set-CellText "E8" "  -1.62%  "
Set-CellText "D9" "0.06324"
Set-CellText "E9" "  +0.86%  "
Set-CellText "E10" "  -2.05%  "
Set-CellText "D11" "0.07793"
Set-CellText "E11" "  +0.52%  "
Set-CellText "D12" "4.500"
Set-CellText "E12" "  +0.40%  "
Set-CellText "D13" "1.651.92"
Set-CellText "E13" "  -0.10%  "
Set-CellText "D14" "1.878.61"
Set-CellText "E14" "  -0.42%  "
Set-CellText "D15" "0.5493"
Set-CellText "E15" "  +0.74%  "
Set-CellText "D16" "0.0₅8201"
Set-CellText "E16" "  +0.81%  "
Set-CellText "D18" "26.079.58"
Set-CellText "E18" "  -0.51%  "
Set-CellText "E19" "  -0.20%  "
Set-CellText "D20" "4.574"
Set-CellText "E20" "  -0.77%  "
Set-CellText "D21" "190.63"
Set-CellText "E21" "  -0.76%  "
Set-CellText "E22" "  -0.21%  "
Set-CellText "E23" "  +0.59%  "
Set-CellText "E24" "  -0.22%  "
Set-CellText "D25" "143.60"
Set-CellText "E25" "  +2.90%  "
Set-CellText "E26" "  +1.51%  "
Set-CellText "D27" "7.226"
Set-CellText "E27" "  -0.67%  "
Set-CellText "D28" "16.03"
Set-CellText "E28" "  -0.54%  "
Set-CellText "D29" "1.428"
Set-CellText "E29" "  -0.86%  "
Set-CellText "D30" "0.05821"
Set-CellText "E30" "  -1.78%  "
Set-CellText "D31" "1.273"
Set-CellText "D32" "3.550"
Set-CellText "E32" "  +0.16%  "
Set-CellText "D33" "3.266"
Set-CellText "E33" "  -0.16%  "
Set-CellText "D34" "1.584"
Set-CellText "E34" "  +0.13%  "
Set-CellText "D35" "0.9469"
Set-CellText "E35" "  -1.46%  "
Set-CellText "B36" "HuobiToken"
Set-CellText "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-CellText "D36" "2.409"
Set-CellText "E36" "  -0.71%  "
Set-CellText "B37" "MXToken"
Set-CellText "C37" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText "D37" "2.777"
Set-CellText "E37" "  +0.29%  "
Set-CellText "D38" "0.5738"
Set-CellText "E38" "  +1.18%  "
Set-CellText "E39" "  +1.18%  "
Set-CellText "B40" "FraxShare"
Set-CellText "C40" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D40" "5.746"
Set-CellText "E40" "  -4.92%  "
Set-CellText "B41" "TrustWalletToken"
Set-CellText "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D41" "0.8423"
Set-CellText "E41" "  -1.63%  "
Set-CellText "E42" "  -0.05%  "
Set-CellText "D43" "103.88"
Set-CellText "E43" "  +3.25%  "
Set-CellText "D44" "1.030.44"
Set-CellText "E44" "  +1.89%  "
Set-CellText "D45" "1.794.84"
Set-CellText "D46" "56.93"
Set-CellText "E46" "  +0.87%  "
Set-CellText "E47" "  +0.01%  "
Set-CellText "D48" "0.4324"
Set-CellText "E48" "  +2.88%  "
Set-CellText "D49" "7.851"
Set-CellText "E49" "  -2.47%  "
Set-CellText "D50" "0.05142"
Set-CellText "E50" "  -0.42%  "
Set-CellText "D51" "1.458"
Set-CellText "E51" "  +1.11%  "
